$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifting existing rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the new row: sedtype_id = -1, sedtype = "Not applicable"
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Update the defined name range to cover the extra row (D64 -> D65)
$n = $wb.Names.Item(1)
$n.RefersTo = "=dbo_sedtype!`$A`$1:`$D`$65"

# Match the resulting selection left behind in the sheet
$ws.Range("A3").Select()
